$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.892.99'
$ws.Range('E2').Value = '  -4.85%  '
$ws.Range('D3').Value = '2.492.51'
$ws.Range('E3').Value = '  -3.25%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '535.51'
$ws.Range('E5').Value = '  -2.78%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '144.17'
$ws.Range('E6').Value = '  -6.89%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.998'
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.574'
$ws.Range('E8').Value = '  -3.19%  '
$ws.Range('D9').Value = '2.517.59'
$ws.Range('E9').Value = '  -2.30%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0999'
$ws.Range('E10').Value = '  -4.07%  '
$ws.Range('E11').Value = '  -2.68%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.48'
$ws.Range('E12').Value = '  -0.41%  '
$ws.Range('E13').Value = '  -3.60%  '
$ws.Range('D14').Value = '2.943.09'
$ws.Range('E14').Value = '  -2.88%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '23.85'
$ws.Range('E15').Value = '  -6.22%  '
$ws.Range('D16').Value = '58.826.79'
$ws.Range('E16').Value = '  -4.82%  '
$ws.Range('E17').Value = '  -3.70%  '
$ws.Range('D18').Value = '2.520.52'
$ws.Range('E18').Value = '  -2.25%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.29'
$ws.Range('E19').Value = '  -2.49%  '
$ws.Range('E20').Value = '  -5.65%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '323.31'
$ws.Range('E21').Value = '  -4.36%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.995'
$ws.Range('E22').Value = '  -0.47%  '
$ws.Range('E23').Value = '  -4.47%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '61.58'
$ws.Range('E24').Value = '  -3.02%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.440'
$ws.Range('E25').Value = '  -10.46%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.997'
$ws.Range('E26').Value = '  -0.22%  '
$ws.Range('B27').Value = 'WrappedeETH'
$ws.Range('C27').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D27').Value = '2.607.78'
$ws.Range('E27').Value = '  -3.26%  '
$ws.Range('B28').Value = 'Kaspa'
$ws.Range('C28').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.161'
$ws.Range('E28').Value = '  -3.71%  '
$ws.Range('E29').Value = '  -4.54%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.82'
$ws.Range('E30').Value = '  -5.75%  '
$ws.Range('D31').Value = '0.0₃0777'
$ws.Range('E31').Value = '  -6.76%  '
$ws.Range('B32').Value = 'PancakeSwap'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.79'
$ws.Range('E32').Value = '  -5.27%  '
$ws.Range('B33').Value = 'Fetch.AI'
$ws.Range('C33').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.23'
$ws.Range('E33').Value = '  -8.36%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '159.78'
$ws.Range('E34').Value = '  -1.78%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.998'
$ws.Range('E35').Value = '  -0.09%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.45'
$ws.Range('E36').Value = '  +2.95%  '
$ws.Range('E37').Value = '  -3.20%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.42'
$ws.Range('E38').Value = '  -9.22%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.62'
$ws.Range('E39').Value = '  -9.09%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.73'
$ws.Range('E40').Value = '  -5.09%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '307.49'
$ws.Range('E41').Value = '  -5.60%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '36.72'
$ws.Range('E42').Value = '  -2.03%  '
$ws.Range('B43').Value = 'SuiNetwork'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.825'
$ws.Range('E43').Value = '  -8.53%  '
$ws.Range('B44').Value = 'Filecoin'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.67'
$ws.Range('E44').Value = '  -6.43%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.993'
$ws.Range('E45').Value = '  -0.47%  '
$ws.Range('B46').Value = 'Mantle'
$ws.Range('C46').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.594'
$ws.Range('E46').Value = '  -1.72%  '
$ws.Range('B47').Value = 'WhiteBITCoin'
$ws.Range('C47').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.77'
$ws.Range('E47').Value = '  -1.56%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '124.93'
$ws.Range('E48').Value = '  +2.30%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0930'
$ws.Range('E49').Value = '  -3.59%  '
$ws.Range('E50').Value = '  -4.69%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0515'
$ws.Range('E51').Value = '  -5.47%  '
